$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.569.25"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "2.610.00"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "321.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.50%  "
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0807"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "3.013.86"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "2.608.26"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.861"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "48.561.17"
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.73"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.136"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0794"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.94"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "2.059.21"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  -4.60%  "
